# Append: 2026-01-26 06:34 JST
# Updates the "ランサーズ" (lancers) listing sheet: refreshes the timestamp,
# replaces the job listings with the newest scrape results, and shrinks the
# sheet from 10 data rows (A2:H11) down to 5 data rows (A2:H6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2026-01-26 06:34:17"

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "製造業向け設備要件定義書の自動生成AIシステムの開発・DB設計支援エンジニア(AI/バックエンド)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5473648"
$ws.Range("G2").Value = 390
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【急募】対話型AI WebアプリMVP開発エンジニア募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5478844"
$ws.Range("G3").Value = 378
$ws.Range("H3").Value = "🔥AI,Ai ◆開発 ◇アプリ"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "自動化システム"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5477084"
$ws.Range("G4").Value = 110
$ws.Range("H4").Value = "◆自動化"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "製造業DXプロダクト開発のプロダクトマネージャー募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5468432"
$ws.Range("G5").Value = 75
$ws.Range("H5").Value = "◆開発"

# --- Row 6 ---------------------------------------------------------------
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "【医療保険】オンライン資格確認・請求端末セットアップ依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5478715"
$ws.Range("G6").Value = 13
$ws.Range("H6").ClearContents()

# --- Drop the now-stale rows 7:11 ----------------------------------------
$ws.Range("A7:H11").EntireRow.Delete()

# --- Rebuild the hyperlinks so only F2:F6 remain --------------------------
# (row deletion above does not automatically drop hyperlinks that pointed at
# the removed rows, so clear them all out and re-create the 5 still needed)
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5473648")
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5478844")
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5477084")
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5468432")
$ws.Range("F5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5478715")
$ws.Range("F6").Style = "Hyperlink"

# --- Column width tweaks ---------------------------------------------------
# Column B: 50 -> 51, Column H: 14 -> 17 (raw OOXML width units). Excel's
# ColumnWidth property is offset from the stored width by ~5/6 of a
# character, so back that out to land on the exact target widths.
$ws.Columns.Item(2).ColumnWidth = 51 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 17 - (5/6)
